# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Reorders the "Periodo Mora" column (E16:E50) from descending (2003 -> 1705)
# to ascending (1705 -> 2003) order, matching the new shared-string ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}
